# "calculo de estatisticas de erro parcial"
# Plan1 holds the calibration time series: column A is the year index
# (originally 0..10, driven by A2 plus a chain of "=prev+1" formulas in
# A3:A12) and column B the observed values used for calibration. Rebase
# the series on real calendar years (2007..2017) by just updating the
# seed cell - the dependent formulas recompute themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("A2").Value = 2007

# Leave the active selection where the analyst ended up after the edit.
$ws.Activate()
$ws.Range("A3").Select()
